$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet / update "through" date references from 10-17 to 10-18
$ws.Name = "Through 2022-10-18"
$ws.Range("A11").Value = "October (through 10-18)"

# Update September I10
$ws.Range("I10").Value = 145

# Update October row (row 11)
$ws.Range("B11").Value = 17
$ws.Range("D11").Value = 31
$ws.Range("E11").Value = 46
$ws.Range("F11").Value = 28
$ws.Range("G11").Value = 86
$ws.Range("H11").Value = 113
$ws.Range("I11").Value = 60

# Update Total row (row 12)
$ws.Range("B12").Value = 243
$ws.Range("D12").Value = 658
$ws.Range("E12").Value = 594
$ws.Range("F12").Value = 450
$ws.Range("G12").Value = 987
$ws.Range("H12").Value = 1360
$ws.Range("I12").Value = 1337
